$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.316.32'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.807.98'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.00'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5151'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3985'
$ws.Range("E8").Value = '  +3.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07839'
$ws.Range("E9").Value = '  -5.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.111'
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.00'
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.344'
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.003'
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("E14").Value = '  -3.21%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.818.37'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.306'
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.58'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001085'
$ws.Range("E18").Value = '  -3.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06561'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.007'
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").Value = '28.346.33'
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -2.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.232'
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.95'
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.52'
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("D28").Value = '2.019.29'
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.413'
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.89'
$ws.Range("E30").Value = '  +1.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1100'
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.064'
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.672'
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.563'
$ws.Range("E34").Value = '  -2.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07185'
$ws.Range("E35").Value = '  -4.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.169'
$ws.Range("E36").Value = '  +4.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02346'
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2180'
$ws.Range("E38").Value = '  -1.81%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.046'
$ws.Range("E39").Value = '  -3.75%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.57'
$ws.Range("E40").Value = '  -4.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6186'
$ws.Range("E41").Value = '  -3.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.154'
$ws.Range("E43").Value = '  -2.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.24'
$ws.Range("E44").Value = '  -2.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5986'
$ws.Range("E45").Value = '  -3.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.305'
$ws.Range("E46").Value = '  -6.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.747'
$ws.Range("E47").Value = '  -1.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.23'
$ws.Range("E48").Value = '  -1.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.220'
$ws.Range("E49").Value = '  +1.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.920'
$ws.Range("E50").Value = '  -4.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06833'
$ws.Range("E51").Value = '  -1.82%  '
